$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4836
$ws.Range("H40").Value = 1991.875
$ws.Range("I40").Value = 1905.8334
$ws.Range("J40").Value = 2250
$ws.Range("K40").Value = 1905.8334
$ws.Range("L40").Value = 2250
$ws.Range("M40").Value = -1730.8334
$ws.Range("N40").Value = -2600
$ws.Range("H64").Value = 4254.706
$ws.Range("I64").Value = 4249.1665
$ws.Range("J64").Value = 4268
$ws.Range("K64").Value = 4249.1665
$ws.Range("L64").Value = 4268
$ws.Range("M64").Value = -4001.1665
$ws.Range("N64").Value = -4764
$ws.Range("H67").Value = 4254.706
$ws.Range("I67").Value = 4249.1665
$ws.Range("J67").Value = 4268
$ws.Range("K67").Value = 4249.1665
$ws.Range("L67").Value = 4268
$ws.Range("M67").Value = -3391.1665
$ws.Range("N67").Value = -5984
$ws.Range("H112").Value = 2733.8096
$ws.Range("I112").Value = 1033.3334
$ws.Range("K112").Value = 3100.0002
$ws.Range("M112").Value = -1992.0002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2142.4688
$ws.Range("I63").Value = 2045.409
$ws.Range("J63").Value = 2356
$ws.Range("K63").Value = 2045.409
$ws.Range("L63").Value = 2356
$ws.Range("M63").Value = -1359.409
$ws.Range("N63").Value = -3728
$ws.Range("H66").Value = 2142.4688
$ws.Range("I66").Value = 2045.409
$ws.Range("J66").Value = 2356
$ws.Range("K66").Value = 10227.045
$ws.Range("L66").Value = 11780
$ws.Range("M66").Value = -6795.045
$ws.Range("N66").Value = -18644
$ws.Range("H74").Value = 1468.1052
$ws.Range("I74").Value = 1080.875
$ws.Range("K74").Value = 1080.875
$ws.Range("M74").Value = -206.875
$ws.Range("H77").Value = 1468.1052
$ws.Range("I77").Value = 1080.875
$ws.Range("K77").Value = 5404.375
$ws.Range("M77").Value = -1036.375
$ws.Range("H114").Value = 24832.666
$ws.Range("I114").Value = 10000
$ws.Range("J114").Value = 27799.2
$ws.Range("K114").Value = 10000
$ws.Range("L114").Value = 27799.2
$ws.Range("M114").Value = -5661
$ws.Range("N114").Value = -36477.2
$ws.Range("H132").Value = 2687.5
$ws.Range("I132").Value = 2205.5
$ws.Range("K132").Value = 6616.5
$ws.Range("M132").Value = -4086.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 22578
$ws.Range("J140").Value = 22578
$ws.Range("L140").Value = 22578
$ws.Range("N140").Value = -32938

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1262.4
$ws.Range("I31").Value = 1105.909
$ws.Range("J31").Value = 3844.5
$ws.Range("K31").Value = 1105.909
$ws.Range("L31").Value = 3844.5
$ws.Range("M31").Value = -810.9090000000001
$ws.Range("N31").Value = -4434.5
$ws.Range("H34").Value = 1262.4
$ws.Range("I34").Value = 1105.909
$ws.Range("J34").Value = 3844.5
$ws.Range("K34").Value = 1105.909
$ws.Range("L34").Value = 3844.5
$ws.Range("M34").Value = -903.9090000000001
$ws.Range("N34").Value = -4248.5
$ws.Range("H99").Value = 1875
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -4996
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H126").Value = 1875
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2780
$ws.Range("N126").Value = -10940
$ws.Range("H134").Value = 31252682
$ws.Range("I134").Value = 2992.5454
$ws.Range("K134").Value = 8977.636200000001
$ws.Range("M134").Value = -6442.636200000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50000120
$ws.Range("I12").Value = 200000210
$ws.Range("J12").Value = 91.53333000000001
$ws.Range("K12").Value = 600000630
$ws.Range("L12").Value = 274.59999
$ws.Range("M12").Value = -600000457
$ws.Range("N12").Value = -620.59999
$ws.Range("H13").Value = 375.8
$ws.Range("I13").Value = 220
$ws.Range("K13").Value = 660
$ws.Range("M13").Value = -492
$ws.Range("H50").Value = 290
$ws.Range("I50").Value = 53.333332
$ws.Range("J50").Value = 1000
$ws.Range("K50").Value = 159.999996
$ws.Range("L50").Value = 3000
$ws.Range("M50").Value = 321.000004
$ws.Range("N50").Value = -3962
$ws.Range("H51").Value = 1214.8572
$ws.Range("I51").Value = 1167.3334
$ws.Range("K51").Value = 3502.0002
$ws.Range("M51").Value = -3042.0002
$ws.Range("H53").Value = 290
$ws.Range("I53").Value = 53.333332
$ws.Range("J53").Value = 1000
$ws.Range("K53").Value = 159.999996
$ws.Range("L53").Value = 3000
$ws.Range("M53").Value = 321.000004
$ws.Range("N53").Value = -3962
$ws.Range("H104").Value = 5232.5
$ws.Range("J104").Value = 5255.4443
$ws.Range("L104").Value = 15766.3329
$ws.Range("N104").Value = -21008.3329
$ws.Range("H107").Value = 14758.429
$ws.Range("J107").Value = 25372.25
$ws.Range("L107").Value = 76116.75
$ws.Range("N107").Value = -79956.75
$ws.Range("H139").Value = 1918.4762
$ws.Range("I139").Value = 2018
$ws.Range("K139").Value = 6054
$ws.Range("M139").Value = -914

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2780.8
$ws.Range("I80").Value = 1733.3334
$ws.Range("J80").Value = 3479.111
$ws.Range("K80").Value = 1733.3334
$ws.Range("L80").Value = 3479.111
$ws.Range("M80").Value = -735.3334
$ws.Range("N80").Value = -5475.111
$ws.Range("H83").Value = 2780.8
$ws.Range("I83").Value = 1733.3334
$ws.Range("J83").Value = 3479.111
$ws.Range("K83").Value = 8666.666999999999
$ws.Range("L83").Value = 17395.555
$ws.Range("M83").Value = -3674.666999999999
$ws.Range("N83").Value = -27379.555
$ws.Range("H103").Value = 34966.668
$ws.Range("J103").Value = 34966.668
$ws.Range("L103").Value = 34966.668
$ws.Range("N103").Value = -37310.668

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1758.2307
$ws.Range("I7").Value = 1594.2
$ws.Range("K7").Value = 1594.2
$ws.Range("M7").Value = -1482.2
$ws.Range("H40").Value = 2593
$ws.Range("I40").Value = 1747.5264
$ws.Range("J40").Value = 10625
$ws.Range("K40").Value = 1747.5264
$ws.Range("L40").Value = 10625
$ws.Range("M40").Value = -1611.5264
$ws.Range("N40").Value = -10897
$ws.Range("H126").Value = 1758.2307
$ws.Range("I126").Value = 1594.2
$ws.Range("K126").Value = 4782.6
$ws.Range("M126").Value = -2312.6
$ws.Range("H132").Value = 125611.78
$ws.Range("I132").Value = 28125.75
$ws.Range("J132").Value = 203600.6
$ws.Range("K132").Value = 84377.25
$ws.Range("L132").Value = 610801.8
$ws.Range("M132").Value = -81847.25
$ws.Range("N132").Value = -615861.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 823800.6
$ws.Range("J4").Value = 823800.6
$ws.Range("L4").Value = 823800.6
$ws.Range("N4").Value = -824026.6
$ws.Range("H126").Value = 58825844
$ws.Range("I126").Value = 76925190
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 230775570
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -230773100
$ws.Range("N126").Value = -13790
$ws.Range("H132").Value = 2397.3333
$ws.Range("I132").Value = 1286.6154
$ws.Range("J132").Value = 3428.7144
$ws.Range("K132").Value = 3859.8462
$ws.Range("L132").Value = 10286.1432
$ws.Range("M132").Value = -1329.8462
$ws.Range("N132").Value = -15346.1432
$ws.Range("H136").Value = 1073.5
$ws.Range("I136").Value = 1079.3334
$ws.Range("K136").Value = 3238.0002
$ws.Range("M136").Value = -688.0001999999999
